$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "309.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.73%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.89%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.072"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.33%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07912"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.46%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.030"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.24%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.402"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.18%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.257"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.17%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.13%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9283"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.12%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1286"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.06%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1889"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.89%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08766"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.38%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03460"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.29%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09757"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.60%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001399"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.19%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006093"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.72%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.589"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.30%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.55%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.11%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.016"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.17%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2521"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.25%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04337"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.88%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.81%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004621"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.87%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "176.86%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.87%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05080"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.58%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007571"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.28%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009937"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.11%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1369"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.49%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002027"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.06%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008830"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.45%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006510"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.11%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.60%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003007"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "8.82%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001205"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.68%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.60%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.60%"
